$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "61.113.79"
$ws.Cells.Item(2,5).Value = "  -1.88%  "

$ws.Cells.Item(3,4).Value = "3.433.18"
$ws.Cells.Item(3,5).Value = "  -1.03%  "

$ws.Cells.Item(4,4).Value = "'0.999"
$ws.Cells.Item(4,5).Value = "  +0.00%  "

$ws.Cells.Item(5,4).Value = "'573.34"
$ws.Cells.Item(5,5).Value = "  -1.36%  "

$ws.Cells.Item(6,4).Value = "'142.88"
$ws.Cells.Item(6,5).Value = "  -4.71%  "

$ws.Cells.Item(7,4).Value = "3.431.58"
$ws.Cells.Item(7,5).Value = "  -1.01%  "

$ws.Cells.Item(8,5).Value = "  +0.01%  "

$ws.Cells.Item(9,4).Value = "'0.478"
$ws.Cells.Item(9,5).Value = "  +0.34%  "

$ws.Cells.Item(10,4).Value = "'7.55"
$ws.Cells.Item(10,5).Value = "  -1.87%  "

$ws.Cells.Item(11,5).Value = "  +0.38%  "

$ws.Cells.Item(12,5).Value = "  -1.30%  "

$ws.Cells.Item(13,4).Value = "4.024.55"
$ws.Cells.Item(13,5).Value = "  -0.92%  "

$ws.Cells.Item(14,4).Value = "'28.30"
$ws.Cells.Item(14,5).Value = "  +2.61%  "

$ws.Cells.Item(15,5).Value = "  -0.16%  "

$ws.Cells.Item(16,2).Value = "WrappedEther"
$ws.Cells.Item(16,3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16,4).Value = "3.455.23"
$ws.Cells.Item(16,5).Value = "  -0.22%  "

$ws.Cells.Item(17,2).Value = "ShibaInu"
$ws.Cells.Item(17,3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Cells.Item(17,4).Value = "'0.0000172"
$ws.Cells.Item(17,5).Value = "  -2.47%  "

$ws.Cells.Item(18,4).Value = "61.250.02"
$ws.Cells.Item(18,5).Value = "  -1.62%  "

$ws.Cells.Item(19,5).Value = "  +2.40%  "

$ws.Cells.Item(20,4).Value = "'14.41"
$ws.Cells.Item(20,5).Value = "  +1.90%  "

$ws.Cells.Item(21,4).Value = "'9.40"
$ws.Cells.Item(21,5).Value = "  -1.95%  "

$ws.Cells.Item(22,4).Value = "'398.10"
$ws.Cells.Item(22,5).Value = "  +2.34%  "

$ws.Cells.Item(23,4).Value = "'0.566"
$ws.Cells.Item(23,5).Value = "  +0.03%  "

$ws.Cells.Item(24,4).Value = "'73.13"
$ws.Cells.Item(24,5).Value = "  +0.64%  "

$ws.Cells.Item(25,4).Value = "'0.995"

$ws.Cells.Item(26,4).Value = "'0.0000123"
$ws.Cells.Item(26,5).Value = "  -2.34%  "

$ws.Cells.Item(27,4).Value = "3.590.81"
$ws.Cells.Item(27,5).Value = "  -0.01%  "

$ws.Cells.Item(28,4).Value = "'0.178"
$ws.Cells.Item(28,5).Value = "  -1.05%  "

$ws.Cells.Item(29,4).Value = "'7.54"
$ws.Cells.Item(29,5).Value = "  -4.02%  "

$ws.Cells.Item(30,4).Value = "'0.999"
$ws.Cells.Item(30,5).Value = "  -0.23%  "

$ws.Cells.Item(31,4).Value = "'8.16"
$ws.Cells.Item(31,5).Value = "  -1.80%  "

$ws.Cells.Item(32,4).Value = "'2.17"
$ws.Cells.Item(32,5).Value = "  -0.62%  "

$ws.Cells.Item(33,4).Value = "'1.44"
$ws.Cells.Item(33,5).Value = "  -9.34%  "

$ws.Cells.Item(34,5).Value = "  -0.12%  "

$ws.Cells.Item(35,4).Value = "'23.92"
$ws.Cells.Item(35,5).Value = "  -0.63%  "

$ws.Cells.Item(36,4).Value = "3.462.60"
$ws.Cells.Item(36,5).Value = "  -0.79%  "

$ws.Cells.Item(37,4).Value = "'7.01"
$ws.Cells.Item(37,5).Value = "  -0.86%  "

$ws.Cells.Item(38,4).Value = "'5.13"
$ws.Cells.Item(38,5).Value = "  -4.13%  "

$ws.Cells.Item(39,5).Value = "  -1.65%  "

$ws.Cells.Item(40,4).Value = "'167.03"
$ws.Cells.Item(40,5).Value = "  +0.12%  "

$ws.Cells.Item(41,2).Value = "Hedera"
$ws.Cells.Item(41,3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(41,4).Value = "'0.0785"
$ws.Cells.Item(41,5).Value = "  -1.53%  "

$ws.Cells.Item(42,2).Value = "EnergySwap"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(42,4).Value = "'27.71"
$ws.Cells.Item(42,5).Value = "  +6.67%  "

$ws.Cells.Item(43,4).Value = "'0.802"
$ws.Cells.Item(43,5).Value = "  +0.59%  "

$ws.Cells.Item(44,2).Value = "FirstDigitalUSD"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(44,4).Value = "'1.01"
$ws.Cells.Item(44,5).Value = "  +0.59%  "

$ws.Cells.Item(45,2).Value = "Filecoin"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(45,4).Value = "'4.51"
$ws.Cells.Item(45,5).Value = "  +0.74%  "

$ws.Cells.Item(46,5).Value = "  -1.65%  "

$ws.Cells.Item(47,4).Value = "2.619.59"
$ws.Cells.Item(47,5).Value = "  -1.26%  "

$ws.Cells.Item(48,5).Value = "  -5.49%  "

$ws.Cells.Item(50,4).Value = "'23.01"
$ws.Cells.Item(50,5).Value = "  -3.60%  "

$ws.Cells.Item(51,4).Value = "'2.39"
$ws.Cells.Item(51,5).Value = "  +1.16%  "
